# Refactored Parser structure. Fixed problems with reading size of classes
# and methods: populate the previously-empty "Number of Lines" metrics for
# classNumberOfLines and methodNumberOfLines.
#
# Note: the source values ("0", "1", "3", "7", ...) are stored as TEXT
# (shared strings) in this workbook, not as numbers - entering them via
# Range.Value directly would auto-coerce numeric-looking strings into the
# Number cell type. To preserve the text type we build the value with a
# text formula ("7") and then collapse the formula down to a static value
# with a values-only paste, which keeps the text type.
$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# --- classNumberOfLines: CustomerRepository had 0 lines -> 1 ---
$wsClass = $wb.Worksheets.Item("classNumberOfLines")
Set-TextValue $wsClass.Cells.Item(4, 2) "1"

# --- methodNumberOfLines ---
$wsMethod = $wb.Worksheets.Item("methodNumberOfLines")

Set-TextValue $wsMethod.Cells.Item(2, 3) "7"    # DevBootstrap.run(String[])
Set-TextValue $wsMethod.Cells.Item(9, 3) "1"    # PaymentApplicationTests()

# Customer getters/setters (toString, builder, getId, getName,
# getAmountAvailable, getAmountReserved, setId, setName,
# setAmountAvailable, setAmountReserved)
for ($r = 10; $r -le 19; $r++) {
    Set-TextValue $wsMethod.Cells.Item($r, 3) "3"
}

Set-TextValue $wsMethod.Cells.Item(20, 3) "2"   # Customer()
Set-TextValue $wsMethod.Cells.Item(21, 3) "6"   # Customer(Long, String, int, int)
Set-TextValue $wsMethod.Cells.Item(33, 3) "1"   # PaymentApplication()
